# Daily COVID-19 data refresh + re-sort of a handful of countries whose
# "Casos totales" (column B) ranking changed relative to their neighbours.
# Country names in column A are rewritten alongside the numbers so that rows
# that swapped rank (Panama/Israel, the Surinam block, Seychelles/Montserrat,
# Papua Nueva Guinea/Islas Virgenes Britanicas) end up with the right label.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Estados Unidos (numbers updated)
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 2089402
$ws.Range("C4").Value = 23001
$ws.Range("D4").Value = 814115
$ws.Range("E4").Value = 1159258
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 899
$ws.Range("H4").Value = 116029

# Row 5: Brasil (numbers updated)
$ws.Range("A5").Value = "Brasil"
$ws.Range("B5").Value = 805649
$ws.Range("C5").Value = 30465
$ws.Range("D5").Value = 396692
$ws.Range("E5").Value = 367899
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 1261
$ws.Range("H5").Value = 41058

# Row 40: Argentina (numbers updated)
$ws.Range("A40").Value = "Argentina"
$ws.Range("B40").Value = 27373
$ws.Range("C40").Value = 1386
$ws.Range("D40").Value = 8332
$ws.Range("E40").Value = 18276
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 30
$ws.Range("H40").Value = 765

# Row 47: Israel -> Panama
$ws.Range("A47").Value = "Panama"
$ws.Range("B47").Value = 18586
$ws.Range("C47").Value = 697
$ws.Range("D47").Value = 10977
$ws.Range("E47").Value = 7191
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 5
$ws.Range("H47").Value = 418

# Row 48: Panama -> Israel
$ws.Range("A48").Value = "Israel"
$ws.Range("B48").Value = 18569
$ws.Range("C48").Value = 214
$ws.Range("D48").Value = 15250
$ws.Range("E48").Value = 3019
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 1
$ws.Range("H48").Value = 300

# Row 93: Venezuela (numbers updated)
$ws.Range("A93").Value = "Venezuela"
$ws.Range("B93").Value = 2814
$ws.Range("C93").Value = 76
$ws.Range("D93").Value = 487
$ws.Range("E93").Value = 2304
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 0
$ws.Range("H93").Value = 23

# Row 153: Libia (numbers updated)
$ws.Range("A153").Value = "Libia"
$ws.Range("B153").Value = 393
$ws.Range("C153").Value = 15
$ws.Range("D153").Value = 59
$ws.Range("E153").Value = 329
$ws.Range("F153").Value = 0
$ws.Range("G153").Value = 0
$ws.Range("H153").Value = 5

# Row 164: Islas Caimanes (numbers updated)
$ws.Range("A164").Value = "Islas Caimanes"
$ws.Range("B164").Value = 186
$ws.Range("C164").Value = 6
$ws.Range("D164").Value = 112
$ws.Range("E164").Value = 73
$ws.Range("F164").Value = 0
$ws.Range("G164").Value = 0
$ws.Range("H164").Value = 1

# Row 166: Siria -> Surinam
$ws.Range("A166").Value = "Surinam"
$ws.Range("B166").Value = 168
$ws.Range("C166").Value = 24
$ws.Range("D166").Value = 9
$ws.Range("E166").Value = 157
$ws.Range("F166").Value = 0
$ws.Range("G166").Value = 0
$ws.Range("H166").Value = 2

# Row 167: Guadalupe -> Siria
$ws.Range("A167").Value = "Siria"
$ws.Range("B167").Value = 164
$ws.Range("C167").Value = 12
$ws.Range("D167").Value = 68
$ws.Range("E167").Value = 90
$ws.Range("F167").Value = 0
$ws.Range("G167").Value = 0
$ws.Range("H167").Value = 6

# Row 168: Comoras -> Guadalupe
$ws.Range("A168").Value = "Guadalupe"
$ws.Range("B168").Value = 164
$ws.Range("C168").Value = 0
$ws.Range("D168").Value = 144
$ws.Range("E168").Value = 6
$ws.Range("F168").Value = 0
$ws.Range("G168").Value = 0
$ws.Range("H168").Value = 14

# Row 169: Guyana -> Comoras
$ws.Range("A169").Value = "Comoras"
$ws.Range("B169").Value = 162
$ws.Range("C169").Value = 0
$ws.Range("D169").Value = 97
$ws.Range("E169").Value = 63
$ws.Range("F169").Value = 0
$ws.Range("G169").Value = 0
$ws.Range("H169").Value = 2

# Row 170: Surinam -> Guyana
$ws.Range("A170").Value = "Guyana"
$ws.Range("B170").Value = 158
$ws.Range("C170").Value = 2
$ws.Range("D170").Value = 92
$ws.Range("E170").Value = 54
$ws.Range("F170").Value = 0
$ws.Range("G170").Value = 0
$ws.Range("H170").Value = 12

# Row 171: Bermudas (numbers updated)
$ws.Range("A171").Value = "Bermudas"
$ws.Range("B171").Value = 141
$ws.Range("C171").Value = 0
$ws.Range("D171").Value = 127
$ws.Range("E171").Value = 5
$ws.Range("F171").Value = 0
$ws.Range("G171").Value = 0
$ws.Range("H171").Value = 9

# Row 176: Bahamas (numbers updated)
$ws.Range("A176").Value = "Bahamas"
$ws.Range("B176").Value = 103
$ws.Range("C176").Value = 0
$ws.Range("D176").Value = 68
$ws.Range("E176").Value = 24
$ws.Range("F176").Value = 0
$ws.Range("G176").Value = 0
$ws.Range("H176").Value = 11

# Row 210: Montserrat -> Seychelles
$ws.Range("A210").Value = "Seychelles"
$ws.Range("B210").Value = 11
$ws.Range("C210").Value = 0
$ws.Range("D210").Value = 11
$ws.Range("E210").Value = 0
$ws.Range("F210").Value = 0
$ws.Range("G210").Value = 0
$ws.Range("H210").Value = 0

# Row 211: Seychelles -> Montserrat
$ws.Range("A211").Value = "Montserrat"
$ws.Range("B211").Value = 11
$ws.Range("C211").Value = 0
$ws.Range("D211").Value = 10
$ws.Range("E211").Value = 0
$ws.Range("F211").Value = 0
$ws.Range("G211").Value = 0
$ws.Range("H211").Value = 1

# Row 213: Islas Virgenes Britanicas -> Papua Nueva Guinea
$ws.Range("A213").Value = "Papua Nueva Guinea"
$ws.Range("B213").Value = 8
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 8
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 0

# Row 214: Papua Nueva Guinea -> Islas Virgenes Britanicas
$ws.Range("A214").Value = "Islas Virgenes Britanicas"
$ws.Range("B214").Value = 8
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 7
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1
